$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $value) {
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

$ws.Range("D2").Value = '42.600.62'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '2.518.52'
$ws.Range("E3").Value = '  -0.98%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '313.77'
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("D6").Value = '99.01'
$ws.Range("E6").Value = '  -2.13%  '
$ws.Range("D7").Value = '0.562'
$ws.Range("E7").Value = '  -1.19%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  -2.51%  '
$ws.Range("D10").Value = '35.13'
$ws.Range("E10").Value = '  -2.46%  '
Set-TextValue $ws.Range("D11") '0.0800'
$ws.Range("E11").Value = '  -0.56%  '
$ws.Range("E12").Value = '  +0.78%  '
$ws.Range("D13").Value = '7.19'
$ws.Range("E13").Value = '  -2.92%  '
$ws.Range("D14").Value = '2.905.79'
$ws.Range("E14").Value = '  -1.00%  '
$ws.Range("D15").Value = '15.24'
$ws.Range("E15").Value = '  -4.02%  '
$ws.Range("D16").Value = '2.473.86'
$ws.Range("E16").Value = '  -5.19%  '
$ws.Range("E17").Value = '  -3.03%  '
$ws.Range("D18").Value = '42.577.54'
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("E19").Value = '  -3.02%  '
$ws.Range("D20").Value = '0.0₃0936'
Set-TextValue $ws.Range("D21") '12.10'
$ws.Range("E21").Value = '  -1.66%  '
Set-TextValue $ws.Range("D22") '68.80'
$ws.Range("E22").Value = '  -0.43%  '
Set-TextValue $ws.Range("D23") '241.10'
$ws.Range("E23").Value = '  -1.43%  '
$ws.Range("E24").Value = '  -2.68%  '
$ws.Range("E25").Value = '  -3.56%  '
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("D27").Value = '25.35'
$ws.Range("E27").Value = '  -3.95%  '
$ws.Range("D28").Value = '2.26'
$ws.Range("E28").Value = '  -4.28%  '
Set-TextValue $ws.Range("D29") '10.00'
$ws.Range("E29").Value = '  -1.04%  '
$ws.Range("D30").Value = '37.77'
$ws.Range("E30").Value = '  -6.52%  '
$ws.Range("D31").Value = '5.84'
$ws.Range("E31").Value = '  +2.44%  '
$ws.Range("D32").Value = '156.54'
$ws.Range("E32").Value = '  -1.41%  '
$ws.Range("D33").Value = '2.69'
$ws.Range("E33").Value = '  -3.38%  '
$ws.Range("E34").Value = '  +0.52%  '
$ws.Range("E35").Value = '  -2.72%  '
$ws.Range("E36").Value = '  -2.28%  '
$ws.Range("E37").Value = '  -5.27%  '
$ws.Range("D38").Value = '17.62'
$ws.Range("E38").Value = '  -2.76%  '
$ws.Range("E39").Value = '  -3.71%  '
$ws.Range("E40").Value = '  -0.61%  '
$ws.Range("D41").Value = '4.18'
$ws.Range("E41").Value = '  +0.72%  '
$ws.Range("D42").Value = '21.54'
$ws.Range("E42").Value = '  -3.37%  '
$ws.Range("E43").Value = '  -0.05%  '
$ws.Range("E44").Value = '  -1.08%  '
$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").Value = '3.22'
$ws.Range("E45").Value = '  -3.03%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '1.997.66'
$ws.Range("E46").Value = '  +1.34%  '
$ws.Range("D47").Value = '8.94'
$ws.Range("E47").Value = '  +0.70%  '
$ws.Range("D48").Value = '2.755.99'
$ws.Range("E48").Value = '  -1.35%  '
$ws.Range("D49").Value = '78.82'
$ws.Range("E49").Value = '  -2.76%  '
$ws.Range("E50").Value = '  -2.79%  '
Set-TextValue $ws.Range("D51") '71.30'
$ws.Range("E51").Value = '  -2.72%  '
